$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value2 = '42.796.07'
$ws.Range('E2').Value2 = '  -0.63%  '
$ws.Range('D3').Value2 = '2.291.79'
$ws.Range('E3').Value2 = '  -0.88%  '
$ws.Range('E4').Value2 = '  -0.05%  '
$ws.Range('D5').Value2 = "'299.87"
$ws.Range('E5').Value2 = '  -0.70%  '
$ws.Range('D6').Value2 = "'96.66"
$ws.Range('E6').Value2 = '  -1.92%  '
$ws.Range('E7').Value2 = '  +0.81%  '
$ws.Range('E8').Value2 = '  -0.04%  '
$ws.Range('E9').Value2 = '  -3.75%  '
$ws.Range('D10').Value2 = "'35.62"
$ws.Range('E10').Value2 = '  -0.12%  '
$ws.Range('E11').Value2 = '  -0.34%  '
$ws.Range('E12').Value2 = '  +0.60%  '
$ws.Range('D13').Value2 = "'17.70"
$ws.Range('E13').Value2 = '  -1.23%  '
$ws.Range('E14').Value2 = '  -2.14%  '
$ws.Range('D15').Value2 = '2.647.72'
$ws.Range('E15').Value2 = '  -0.91%  '
$ws.Range('D16').Value2 = '2.293.28'
$ws.Range('E16').Value2 = '  +1.02%  '
$ws.Range('D18').Value2 = '42.729.64'
$ws.Range('E18').Value2 = '  -0.61%  '
$ws.Range('D19').Value2 = "'12.75"
$ws.Range('E19').Value2 = '  -4.92%  '
$ws.Range('E20').Value2 = '  -0.49%  '
$ws.Range('E21').Value2 = '  -2.40%  '
$ws.Range('D22').Value2 = "'67.76"
$ws.Range('E22').Value2 = '  -0.46%  '
$ws.Range('E23').Value2 = '  -0.02%  '
$ws.Range('D24').Value2 = "'2.12"
$ws.Range('E24').Value2 = '  -1.49%  '
$ws.Range('E25').Value2 = '  +0.02%  '
$ws.Range('E26').Value2 = '  -0.99%  '
$ws.Range('D27').Value2 = "'4.01"
$ws.Range('E27').Value2 = '  -0.42%  '
$ws.Range('D28').Value2 = "'25.18"
$ws.Range('E28').Value2 = '  +0.50%  '
$ws.Range('D29').Value2 = "'166.13"
$ws.Range('E29').Value2 = '  -1.69%  '
$ws.Range('E30').Value2 = '  -1.10%  '
$ws.Range('E31').Value2 = '  -1.64%  '
$ws.Range('D32').Value2 = "'32.86"
$ws.Range('E32').Value2 = '  -1.38%  '
$ws.Range('E33').Value2 = '  +0.10%  '
$ws.Range('D34').Value2 = "'4.81"
$ws.Range('E34').Value2 = '  -1.76%  '
$ws.Range('D35').Value2 = "'5.00"
$ws.Range('E35').Value2 = '  -3.79%  '
$ws.Range('D36').Value2 = "'16.99"
$ws.Range('E36').Value2 = '  -6.88%  '
$ws.Range('E37').Value2 = '  -1.84%  '
$ws.Range('D38').Value2 = "'0.0683"
$ws.Range('E38').Value2 = '  -1.44%  '
$ws.Range('E39').Value2 = '  -1.47%  '
$ws.Range('E40').Value2 = '  -3.55%  '
$ws.Range('E41').Value2 = '  +0.30%  '
$ws.Range('E42').Value2 = '  -1.06%  '
$ws.Range('D43').Value2 = '2.015.29'
$ws.Range('E43').Value2 = '  +1.17%  '
$ws.Range('E44').Value2 = '  -2.62%  '
$ws.Range('D45').Value2 = "'10.06"
$ws.Range('E45').Value2 = '  -0.65%  '
$ws.Range('D46').Value2 = "'2.08"
$ws.Range('E46').Value2 = '  +0.65%  '
$ws.Range('D47').Value2 = "'17.20"
$ws.Range('E47').Value2 = '  -1.85%  '
$ws.Range('E48').Value2 = '  -2.06%  '
$ws.Range('D49').Value2 = "'2.93"
$ws.Range('E49').Value2 = '  -2.72%  '
$ws.Range('D50').Value2 = '2.515.87'
$ws.Range('E50').Value2 = '  -0.83%  '
$ws.Range('D51').Value2 = "'52.95"
$ws.Range('E51').Value2 = '  -3.05%  '
